$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit rotates the fungi observation records currently on rows 5, 6 and 7
# one step down (row 5 -> row 6, row 6 -> row 7, row 7 -> row 5); the row
# numbers/layout stay put, only the per-record field values move:
#   row 5 <- old row 7,  row 6 <- old row 5,  row 7 <- old row 6

# --- Row 5 now holds the record that used to be on row 7 ---
$ws.Range("A5").Value = 111817582
$ws.Range("B5").Value = 88283
$ws.Range("D5").Value = 'NT'
$ws.Range("E5").Value = 655
$ws.Range("F5").Value = 'Oxtungssvamp'
$ws.Range("G5").Value = 'Fistulina hepatica'
$ws.Range("H5").Value = '(Schaeff.) With., nom sanct.'
$ws.Range("I5").NumberFormat = "@"
$ws.Range("I5").Value = '2'
$ws.Range("I5").Style = "Normal"
$ws.Range("Q5").Value = 578498.8708077573
$ws.Range("R5").Value = 6398730.978152275
$ws.Range("AC5").ClearContents() | Out-Null

# --- Row 6 now holds the record that used to be on row 5 ---
$ws.Range("A6").Value = 111817654
$ws.Range("B6").Value = 89363
$ws.Range("D6").Value = 'NT'
$ws.Range("E6").Value = 5445
$ws.Range("F6").Value = 'Ekticka'
$ws.Range("G6").Value = 'Fomitiporia robusta'
$ws.Range("H6").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("I6").NumberFormat = "@"
$ws.Range("I6").Value = '5'
$ws.Range("I6").Style = "Normal"
$ws.Range("Q6").Value = 578449.7681887138
$ws.Range("R6").Value = 6398641.347911141
$ws.Range("AC6").Value = 'På relativt tunn ek.'

# --- Row 7 now holds the record that used to be on row 6 ---
$ws.Range("A7").Value = 111817611
$ws.Range("B7").Value = 89416
$ws.Range("D7").Value = 'LC'
$ws.Range("E7").Value = 1205
$ws.Range("F7").Value = 'Stor aspticka'
$ws.Range("G7").Value = 'Phellinus populicola'
$ws.Range("H7").Value = 'Niemelä'
$ws.Range("I7").NumberFormat = "@"
$ws.Range("I7").Value = '1'
$ws.Range("I7").Style = "Normal"
$ws.Range("Q7").Value = 578480.2128223784
$ws.Range("R7").Value = 6398699.632505047
$ws.Range("AC7").Value = 'Relativt tunn asp.'
